$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.518.36'
$ws.Range("E2").Value = '  +1.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.009.75'
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.05'
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.34'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.51'
$ws.Range("E9").Value = '  -1.80%  '

$ws.Range("E10").Value = '  +1.41%  '

$ws.Range("E11").Value = '  +3.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.521.12'
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("E13").Value = '  +0.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.41'
$ws.Range("E14").Value = '  +3.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000164'
$ws.Range("E15").Value = '  +5.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.505.65'
$ws.Range("E16").Value = '  +1.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.20'
$ws.Range("E17").Value = '  +6.58%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.005.44'
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.80'
$ws.Range("E19").Value = '  +2.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.96'
$ws.Range("E20").Value = '  +1.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.46'
$ws.Range("E21").Value = '  +0.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.498'
$ws.Range("E23").Value = '  +3.77%  '

$ws.Range("E24").Value = '  +3.11%  '

$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0920'
$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.80'
$ws.Range("E28").Value = '  +3.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("E29").Value = '  +5.05%  '

$ws.Range("E30").Value = '  +2.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("E31").Value = '  -5.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.61'
$ws.Range("E32").Value = '  -0.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.73'
$ws.Range("E33").Value = '  +4.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.83'
$ws.Range("E34").Value = '  -1.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.86'

$ws.Range("E36").Value = '  +0.83%  '

$ws.Range("E37").Value = '  +1.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.41'
$ws.Range("E38").Value = '  +1.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.040.27'
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.33'
$ws.Range("E40").Value = '  +1.76%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.85'
$ws.Range("E41").Value = '  +6.23%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.12%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.649'
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.274.14'
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.985'
$ws.Range("E46").Value = '  -1.61%  '

$ws.Range("E47").Value = '  +4.04%  '

$ws.Range("E48").Value = '  +1.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.41'
$ws.Range("E49").Value = '  +1.50%  '

$ws.Range("E50").Value = '  -7.09%  '

$ws.Range("E51").Value = '  +2.46%  '

